# Append the new "obs 54", "MCAR 34,35,36" and "MCAR 60" sections to the
# end of the document, after the existing "mCAR 62" / "Allow the code to
# used in another monitoring report." content.

$d = $word.ActiveDocument

$cr = [string][char]13
$enDash = [string][char]0x2013
$lsq = [string][char]0x2018
$rsq = [string][char]0x2019

# Build the full block of new paragraphs, separated by paragraph marks.
# Each element below becomes its own paragraph.
$paras = @(
    "",
    ("obs 54 " + $enDash + " code change"),
    "Changed the ARefor survey data to use aa boot rather than the hard wired areas",
    "",
    ("MCAR 34,35,36 " + $enDash + " no code change"),
    "Change the volumes used for the Softwood and Hardwood plantations to use published data.",
    "This required changes to the fiji_frl_input.RData",
    "This changed the hwsw_volharv, sw_hvol_parea data.",
    "",
    ("Volumes used in MR have also changed and will need to be changed in the " + $lsq + "defaults" + $rsq + " and web interface to use the published data."),
    ("MCAR 60 " + $enDash + " no code change"),
    "Risk buffer deduction changed from 16% to 21%",
    ("This only affects the entered values " + $enDash + " need to change the " + $lsq + "defaults" + $rsq),
    ""
)

$text = ($paras -join $cr)

# Insert right at the very end of the document body (before the final
# section mark), so the new content lands after the last existing
# paragraph ("Allow the code to used in another monitoring report.").
$insertAt = $d.Content.End - 1
$rng = $d.Range($insertAt, $insertAt)
$rng.InsertAfter($cr + $text)

# The paragraph that previously was last ("Allow the code...") is still
# at its old index; the freshly inserted paragraphs now follow it.
$total = $d.Paragraphs.Count
$startIndex = $total - $paras.Count + 1

for ($i = 0; $i -lt $paras.Count; $i++) {
    $p = $d.Paragraphs.Item($startIndex + $i)
    if ($paras[$i] -eq ("obs 54 " + $enDash + " code change") -or
        $paras[$i] -eq ("MCAR 34,35,36 " + $enDash + " no code change") -or
        $paras[$i] -eq ("MCAR 60 " + $enDash + " no code change")) {
        $p.Range.Style = "Heading 2"
    }
}

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
